# Update the workbook window view geometry (position/size) to match
# the state captured when the workbook was last saved.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$win = $excel.ActiveWindow
$win.Left = 28680
$win.Top = -120
$win.Width = 29040
$win.Height = 15720

# Append the additional simulation result rows (rows 6-13) below the
# existing data (rows 2-5), reusing the same column layout (B:W).
# Row 6
$ws.Range("B6").Value = 543.34176935416
$ws.Range("C6").Value = 35
$ws.Range("D6").Value = 289296
$ws.Range("E6").Value = 447.2035863348277
$ws.Range("F6").Value = 86788.8
$ws.Range("G6").Value = 72324
$ws.Range("H6").Value = 127
$ws.Range("I6").Value = 98
$ws.Range("J6").Value = 41
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 0
$ws.Range("N6").Value = 1233
$ws.Range("O6").Value = 0.29
$ws.Range("P6").Value = 2.4
$ws.Range("Q6").Value = 0.015
$ws.Range("R6").Value = 64000
$ws.Range("S6").Value = 450.8
$ws.Range("T6").Value = 412
$ws.Range("U6").Value = 245
$ws.Range("V6").Value = 170000
$ws.Range("W6").Value = 77000

# Row 7
$ws.Range("B7").Value = 642.5622171730914
$ws.Range("C7").Value = 77
$ws.Range("D7").Value = 342216
$ws.Range("E7").Value = 447.9235806389746
$ws.Range("F7").Value = 78321.60000000001
$ws.Range("G7").Value = 85554.00000000001
$ws.Range("H7").Value = 111
$ws.Range("I7").Value = 49
$ws.Range("J7").Value = 30
$ws.Range("K7").Value = 127
$ws.Range("L7").Value = 49
$ws.Range("M7").Value = 28
$ws.Range("N7").Value = 1233
$ws.Range("O7").Value = 0.29
$ws.Range("P7").Value = 2.4
$ws.Range("Q7").Value = 0.015
$ws.Range("R7").Value = 64000
$ws.Range("S7").Value = 450.8
$ws.Range("T7").Value = 412
$ws.Range("U7").Value = 245
$ws.Range("V7").Value = 170000
$ws.Range("W7").Value = 77000

# Row 8
$ws.Range("B8").Value = 495.5216781932683
$ws.Range("C8").Value = 63
$ws.Range("D8").Value = 347913.9
$ws.Range("E8").Value = 415.7809128630705
$ws.Range("F8").Value = 87149.5
$ws.Range("G8").Value = 65191.10000000001
$ws.Range("H8").Value = 112
$ws.Range("I8").Value = 49
$ws.Range("J8").Value = 5
$ws.Range("K8").Value = 39
$ws.Range("L8").Value = 47
$ws.Range("M8").Value = 4
$ws.Range("N8").Value = 1233
$ws.Range("O8").Value = 0.29
$ws.Range("P8").Value = 2.4
$ws.Range("Q8").Value = 0.015
$ws.Range("R8").Value = 64000
$ws.Range("S8").Value = 450.8
$ws.Range("T8").Value = 412
$ws.Range("U8").Value = 245
$ws.Range("V8").Value = 170000
$ws.Range("W8").Value = 77000

# Row 9
$ws.Range("B9").Value = 545.4168935434984
$ws.Range("C9").Value = 63
$ws.Range("D9").Value = 599686.5
$ws.Range("E9").Value = 441.1219917012447
$ws.Range("F9").Value = 96187
$ws.Range("G9").Value = 72436.70000000001
$ws.Range("H9").Value = 111
$ws.Range("I9").Value = 49
$ws.Range("J9").Value = 24
$ws.Range("K9").Value = 142
$ws.Range("L9").Value = 49
$ws.Range("M9").Value = 20
$ws.Range("N9").Value = 1233
$ws.Range("O9").Value = 0.29
$ws.Range("P9").Value = 2.4
$ws.Range("Q9").Value = 0.015
$ws.Range("R9").Value = 64000
$ws.Range("S9").Value = 450.8
$ws.Range("T9").Value = 412
$ws.Range("U9").Value = 245
$ws.Range("V9").Value = 170000
$ws.Range("W9").Value = 77000

# Row 10
$ws.Range("B10").Value = 495.5216781932683
$ws.Range("C10").Value = 63
$ws.Range("D10").Value = 347913.9
$ws.Range("E10").Value = 415.7809128630705
$ws.Range("F10").Value = 87149.5
$ws.Range("G10").Value = 65191.10000000001
$ws.Range("H10").Value = 112
$ws.Range("I10").Value = 49
$ws.Range("J10").Value = 5
$ws.Range("K10").Value = 39
$ws.Range("L10").Value = 47
$ws.Range("M10").Value = 4
$ws.Range("N10").Value = 1233
$ws.Range("O10").Value = 0.29
$ws.Range("P10").Value = 2.4
$ws.Range("Q10").Value = 0.015
$ws.Range("R10").Value = 64000
$ws.Range("S10").Value = 450.8
$ws.Range("T10").Value = 412
$ws.Range("U10").Value = 245
$ws.Range("V10").Value = 170000
$ws.Range("W10").Value = 77000

# Row 11
$ws.Range("B11").Value = 514.5199814358043
$ws.Range("C11").Value = 50.4
$ws.Range("D11").Value = 295132.8
$ws.Range("E11").Value = 442.0016597510373
$ws.Range("F11").Value = 115485.6
$ws.Range("G11").Value = 68355.60000000002
$ws.Range("H11").Value = 111
$ws.Range("I11").Value = 50
$ws.Range("J11").Value = 25
$ws.Range("K11").Value = 96
$ws.Range("L11").Value = 48
$ws.Range("M11").Value = 17
$ws.Range("N11").Value = 1233
$ws.Range("O11").Value = 0.29
$ws.Range("P11").Value = 2.4
$ws.Range("Q11").Value = 0.015
$ws.Range("R11").Value = 45000
$ws.Range("S11").Value = 450.8
$ws.Range("T11").Value = 412
$ws.Range("U11").Value = 245
$ws.Range("V11").Value = 170000
$ws.Range("W11").Value = 100000

# Row 12
$ws.Range("B12").Value = 746.3673234890744
$ws.Range("C12").Value = 126
$ws.Range("D12").Value = 260632.08
$ws.Range("E12").Value = 449.5668049792531
$ws.Range("F12").Value = 72841.68000000001
$ws.Range("G12").Value = 99435.60000000001
$ws.Range("H12").Value = 111
$ws.Range("I12").Value = 42
$ws.Range("J12").Value = 27
$ws.Range("K12").Value = 320
$ws.Range("L12").Value = 54
$ws.Range("M12").Value = 75
$ws.Range("N12").Value = 1233
$ws.Range("O12").Value = 0.29
$ws.Range("P12").Value = 2.4
$ws.Range("Q12").Value = 0.015
$ws.Range("R12").Value = 90000
$ws.Range("S12").Value = 450.8
$ws.Range("T12").Value = 412
$ws.Range("U12").Value = 245
$ws.Range("V12").Value = 170000
$ws.Range("W12").Value = 65000

# Row 13
$ws.Range("B13").Value = 458.1427332831215
$ws.Range("C13").Value = 46.66666666666666
$ws.Range("D13").Value = 192384
$ws.Range("E13").Value = 446.442410462481
$ws.Range("F13").Value = 141631.2
$ws.Range("G13").Value = 60966
$ws.Range("H13").Value = 164
$ws.Range("I13").Value = 45
$ws.Range("J13").Value = 4
$ws.Range("K13").Value = 127
$ws.Range("L13").Value = 53
$ws.Range("M13").Value = 29
$ws.Range("N13").Value = 1233
$ws.Range("O13").Value = 0.29
$ws.Range("P13").Value = 2.4
$ws.Range("Q13").Value = 0.015
$ws.Range("R13").Value = 50000
$ws.Range("S13").Value = 450.8
$ws.Range("T13").Value = 412
$ws.Range("U13").Value = 245
$ws.Range("V13").Value = 170000
$ws.Range("W13").Value = 120000

# Restore the active cell / selection on the worksheet to the cell that
# was selected when the workbook was saved.
[void]$ws.Range("M16").Select()
